$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.347.28"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "3.950.47"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'470.94"
$ws.Range("E5").Value = "  +8.24%  "
$ws.Range("D6").Value = "'146.24"
$ws.Range("E6").Value = "  +5.08%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.734"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  +10.26%  "
$ws.Range("D11").Value = "'0.0000347"
$ws.Range("E11").Value = "  +10.25%  "
$ws.Range("D12").Value = "'43.42"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "4.573.13"
$ws.Range("E13").Value = "  +4.52%  "
$ws.Range("D14").Value = "'10.44"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "3.962.98"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'19.87"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").Value = "67.537.16"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "'434.50"
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("D22").Value = "'3.40"
$ws.Range("E22").Value = "  +4.95%  "
$ws.Range("D23").Value = "'14.47"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "'87.52"
$ws.Range("E24").Value = "  +2.18%  "
$ws.Range("D25").Value = "'3.61"
$ws.Range("E25").Value = "  +6.90%  "
$ws.Range("E26").Value = "  +4.91%  "
$ws.Range("D27").Value = "'10.27"
$ws.Range("E27").Value = "  +5.06%  "
$ws.Range("D28").Value = "'9.86"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'720.54"
$ws.Range("E29").Value = "  -0.92%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "'13.51"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").Value = "'2.75"
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "'42.57"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0844"
$ws.Range("E34").Value = "  +25.97%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'58.03"
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.151"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'5.36"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").Value = "'0.0477"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.142"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").Value = "'3.54"
$ws.Range("E42").Value = "  +7.01%  "
$ws.Range("D43").Value = "'0.338"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'2.82"
$ws.Range("E45").Value = "  +5.29%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.54"
$ws.Range("E46").Value = "  -8.16%  "
$ws.Range("D47").Value = "'2.20"
$ws.Range("E47").Value = "  +6.12%  "
$ws.Range("D48").Value = "'149.07"
$ws.Range("E48").Value = "  +5.08%  "
$ws.Range("E49").Value = "  -3.26%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "'25.83"
$ws.Range("E51").Value = "  +4.39%  "
